$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5: Registro "104910020" was not found in the source -> mark as "Não encontrado"
$ws.Range("E5").Value = "Não encontrado"

# Row 5: since the registry record could not be located, the PDF status is now "Pendente"
$ws.Range("F5").Value = "Pendente"
